# TP4 - Final arternativo Nro5 - El despertar del Grafo
#
# Add a totals row under the "1000y75" statistics table that sums the
# occurrence counts for each coloring algorithm (Matula, Secuencial
# aleatorio, Welsh-Powell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1000y75")

$ws.Range("B20").Formula = "=SUM(B4:B19)"
$ws.Range("C20").Formula = "=SUM(C4:C19)"
$ws.Range("D20").Formula = "=SUM(D4:D19)"

$ws.Activate() | Out-Null
$ws.Range("F8").Select() | Out-Null
